$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 41: Mirko - "ultimi ritocchi alla gerarchia, passaggio ai cpp"
$ws.Cells.Item(41, 1).Value = 43133
$ws.Cells.Item(41, 2).Value = "Mirko"
$ws.Cells.Item(41, 3).Value = "ultimi ritocchi alla gerarchia, passaggio ai cpp"
$ws.Cells.Item(41, 4).Value = 0.33333333333333331
$ws.Rows.Item(41).RowHeight = 28.8

# Row 42: Giovanni - "ultimi ritocchi alla gerarchia, passaggio ai cpp, prima implementazione di java"
$ws.Cells.Item(42, 1).Value = 43133
$ws.Cells.Item(42, 2).Value = "Giovanni"
$ws.Cells.Item(42, 3).Value = "ultimi ritocchi alla gerarchia, passaggio ai cpp, prima implementazione di java"
$ws.Cells.Item(42, 4).Value = 0.41666666666666669
$ws.Rows.Item(42).RowHeight = 57.6

# Scroll/selection: move to the newly-entered last cell (clears the old topLeftCell/scroll
# position and matches the new active selection recorded in the workbook view)
[void]$ws.Range("D42").Select()
